$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell text updates for refreshed crypto price/volume data.
# NumberFormat is forced to "@" (Text) before assignment so numeric-looking
# strings (e.g. "21.50", "1.320", "0.02040") are preserved verbatim instead of
# being coerced to numbers by Excel; Style is reset to "Normal" afterwards so no
# residual formatting/style is introduced on the cell.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.536.39'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.869.38'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.31%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '324.71'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4544'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3843'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -1.69%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07829'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.9906'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '21.50'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.43%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.881.53'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.917'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.639'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06935'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('B16').Style = "Normal"
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C16').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '86.55'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -2.10%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('B17').Style = "Normal"
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('C17').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.006'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000009948'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.68'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.40%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.004'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '28.526.64'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.257'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.069'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.102.57'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '153.89'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.12'
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.683'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '117.27'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -1.47%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.882'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -5.45%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09275'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9056'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -3.11%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.275'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.320'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.257'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.63%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.05664'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.149'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02040'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.46%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.612'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -4.03%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.5564'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.49%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1766'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.51%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '9.637'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.07144'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '11.56'
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5249'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -1.05%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.143'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -1.91%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.802'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '111.67'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.436'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +3.89%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.09%  '
$ws.Range('E51').Style = "Normal"
